$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 118
$ws.Cells.Item(3, 6).Value = 1280
$ws.Cells.Item(4, 6).Value = 921
$ws.Cells.Item(5, 6).Value = 963
$ws.Cells.Item(6, 6).Value = 1712
$ws.Cells.Item(7, 6).Value = 376
$ws.Cells.Item(8, 6).Value = 1142
$ws.Cells.Item(9, 6).Value = 47
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(11, 6).Value = 105
$ws.Cells.Item(12, 6).Value = 259
$ws.Cells.Item(13, 6).Value = 34
$ws.Cells.Item(14, 6).Value = 78
$ws.Cells.Item(15, 6).Value = 638
$ws.Cells.Item(16, 6).Value = 131
$ws.Cells.Item(17, 6).Value = 88
$ws.Cells.Item(18, 6).Value = 24
$ws.Cells.Item(21, 6).Value = 100
$ws.Cells.Item(22, 6).Value = 645
$ws.Cells.Item(23, 6).Value = 14
$ws.Cells.Item(24, 6).Value = 625
$ws.Cells.Item(25, 6).Value = 131
$ws.Cells.Item(26, 6).Value = 31
$ws.Cells.Item(27, 6).Value = 836
$ws.Cells.Item(29, 6).Value = 111
$ws.Cells.Item(31, 6).Value = 247
$ws.Cells.Item(33, 6).Value = 7

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 24
$ws.Cells.Item(7, 6).Value = 239
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(11, 6).Value = 113

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 118
$ws.Cells.Item(4, 6).Value = 1280
$ws.Cells.Item(5, 6).Value = 921
$ws.Cells.Item(6, 6).Value = 963
$ws.Cells.Item(7, 6).Value = 1712
$ws.Cells.Item(8, 6).Value = 376
$ws.Cells.Item(9, 6).Value = 1142
$ws.Cells.Item(10, 6).Value = 47
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(13, 6).Value = 105
$ws.Cells.Item(14, 6).Value = 259
$ws.Cells.Item(15, 6).Value = 34
$ws.Cells.Item(16, 6).Value = 78
$ws.Cells.Item(17, 6).Value = 638
$ws.Cells.Item(18, 6).Value = 131
$ws.Cells.Item(19, 6).Value = 88
$ws.Cells.Item(21, 6).Value = 24
$ws.Cells.Item(26, 6).Value = 24
$ws.Cells.Item(27, 6).Value = 239
$ws.Cells.Item(28, 6).Value = 239
$ws.Cells.Item(29, 6).Value = 100
$ws.Cells.Item(30, 6).Value = 645
$ws.Cells.Item(31, 6).Value = 14
$ws.Cells.Item(32, 6).Value = 625
$ws.Cells.Item(33, 6).Value = 131
$ws.Cells.Item(34, 6).Value = 31
$ws.Cells.Item(35, 6).Value = 836
$ws.Cells.Item(38, 6).Value = 2
$ws.Cells.Item(39, 6).Value = 111
$ws.Cells.Item(41, 6).Value = 247
$ws.Cells.Item(43, 6).Value = 113
$ws.Cells.Item(44, 6).Value = 113
$ws.Cells.Item(46, 6).Value = 7

